# feat: add 2022-Q1 data
#
# - Insert a new detail sheet "2022-Q1" (a copy of the "2021-Q4" detail
#   sheet, same fund/columns, new metrics) positioned between "2021-Q4"
#   and "总计".
# - Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#   existing "2021-Q4" row down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Clone the "2021-Q4" detail sheet into a new "2022-Q1" sheet, right
#    after "2021-Q4" (i.e. before "总计").
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Fund code / name (B2, C2) stay the same (513360, same fund) - only the
# position metrics change. Force D2/E2/F2/G2 to remain TEXT (matching the
# source sheet's inlineStr cells) instead of Excel's automatic numeric
# coercion: apply a text number format before the write, then strip the
# resulting style back off so the cells keep the default (unstyled) look.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "6.05"
$newSheet.Range("E2").Value = "99.49"
$newSheet.Range("F2").Value = "5.60"
$newSheet.Range("G2").Value = "0.3388"
$newSheet.Range("D2:G2").ClearFormats()

$newSheet.Range("H2").Value = 6

# ---------------------------------------------------------------------
# 2. "总计" summary sheet: insert a new row 2 for "2022-Q1" above the
#    existing "2021-Q4" row (which shifts down to row 3).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# Match the index-column style (s="2") used by every other data row by
# copying the format from the row that just shifted down, then clear the
# row-insert's inherited header formatting off the remaining cells.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.34

$total.Range("A3").Value = 1

# Restore the original active sheet (the diff leaves <bookViews> untouched,
# so the workbook should still open on "2021-Q4").
$src.Activate()
